# Implement Nicolas' comments on the EN output (05_dimension_overview):
#  - refresh all the computed statistics in the "Strata" table with the
#    latest figures
#  - rename the "ocap"/"idp" strata to "non_pdi"/"pdi"
#  - the table now only needs 9 data rows (it previously had 11), so the
#    two trailing rows are dropped once their figures are folded into the
#    rows above them.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the last two (now redundant) data rows; this also shrinks the
# sheet's used range from A1:O12 down to A1:O10.
$ws.Rows("11:12").Delete()

# Row 2: TOTAL (5-17 y.o.)
$ws.Range("A2").Value = "TOTAL (5-17 y.o.)"
$ws.Range("B2").Value = "All population groups"
$ws.Range("C2").Value = 68663199
$ws.Range("D2").Value = 33
$ws.Range("E2").Value = 22642378
$ws.Range("F2").Value = 0.1
$ws.Range("G2").Value = 75499
$ws.Range("H2").Value = 0.4
$ws.Range("I2").Value = 286415
$ws.Range("J2").Value = 0.9
$ws.Range("K2").Value = 598007
$ws.Range("L2").Value = 65.6
$ws.Range("M2").Value = 45060900
$ws.Range("N2").Value = 34.4
$ws.Range("O2").Value = 23602299

# Row 3: non_pdi (5-17 y.o.)
$ws.Range("A3").Value = "non_pdi (5-17 y.o.)"
$ws.Range("B3").Value = "non_pdi"
$ws.Range("C3").Value = 48523504
$ws.Range("D3").Value = 24.4
$ws.Range("E3").Value = 11822828
$ws.Range("F3").Value = 0.1
$ws.Range("G3").Value = 52304
$ws.Range("H3").Value = 0.5
$ws.Range("I3").Value = 227095
$ws.Range("J3").Value = 0.9
$ws.Range("K3").Value = 423026
$ws.Range("L3").Value = 74.2
$ws.Range("M3").Value = 35998252
$ws.Range("N3").Value = 25.8
$ws.Range("O3").Value = 12525253

# Row 4: pdi (5-17 y.o.)
$ws.Range("A4").Value = "pdi (5-17 y.o.)"
$ws.Range("B4").Value = "pdi"
$ws.Range("C4").Value = 20139695
$ws.Range("D4").Value = 53.7
$ws.Range("E4").Value = 10819550
$ws.Range("F4").Value = 0.1
$ws.Range("G4").Value = 23195
$ws.Range("H4").Value = 0.3
$ws.Range("I4").Value = 59320
$ws.Range("J4").Value = 0.9
$ws.Range("K4").Value = 174981
$ws.Range("L4").Value = 45
$ws.Range("M4").Value = 9062649
$ws.Range("N4").Value = 55
$ws.Range("O4").Value = 11077047

# Row 5: Girls (5-17 y.o.)
$ws.Range("A5").Value = "Girls (5-17 y.o.)"
$ws.Range("B5").Value = "All population groups"
$ws.Range("C5").Value = 30193399
$ws.Range("D5").Value = 34.1
$ws.Range("E5").Value = 10282418
$ws.Range("F5").Value = 0.1
$ws.Range("G5").Value = 15539
$ws.Range("H5").Value = 0.4
$ws.Range("I5").Value = 106744
$ws.Range("J5").Value = 0.8
$ws.Range("K5").Value = 243257
$ws.Range("L5").Value = 64.7
$ws.Range("M5").Value = 19545441
$ws.Range("N5").Value = 35.3
$ws.Range("O5").Value = 10647958

# Row 6: Boys (5-17 y.o.)
$ws.Range("A6").Value = "Boys (5-17 y.o.)"
$ws.Range("B6").Value = "All population groups"
$ws.Range("C6").Value = 38469800
$ws.Range("D6").Value = 31.8
$ws.Range("E6").Value = 12247982
$ws.Range("F6").Value = 0.2
$ws.Range("G6").Value = 62081
$ws.Range("H6").Value = 0.5
$ws.Range("I6").Value = 190004
$ws.Range("J6").Value = 0.9
$ws.Range("K6").Value = 355336
$ws.Range("L6").Value = 66.6
$ws.Range("M6").Value = 25614397
$ws.Range("N6").Value = 33.4
$ws.Range("O6").Value = 12855403

# Row 7: ECE (5 y.o.)
$ws.Range("A7").Value = "ECE (5 y.o.)"
$ws.Range("B7").Value = "All population groups"
$ws.Range("C7").Value = 6978393
$ws.Range("D7").Value = 63.3
$ws.Range("E7").Value = 4420486
$ws.Range("F7").Value = 0.1
$ws.Range("G7").Value = 5418
$ws.Range("H7").Value = 0.4
$ws.Range("I7").Value = 29389
$ws.Range("J7").Value = 0.5
$ws.Range("K7").Value = 32932
$ws.Range("L7").Value = 35.7
$ws.Range("M7").Value = 2490167
$ws.Range("N7").Value = 64.3
$ws.Range("O7").Value = 4488226

# Row 8: Primary school
$ws.Range("A8").Value = "Primary school"
$ws.Range("B8").Value = "All population groups"
$ws.Range("C8").Value = 26408923
$ws.Range("D8").Value = 22.4
$ws.Range("E8").Value = 5914019
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 10888
$ws.Range("H8").Value = 0.4
$ws.Range("I8").Value = 97051
$ws.Range("J8").Value = 1.1
$ws.Range("K8").Value = 299759
$ws.Range("L8").Value = 76.1
$ws.Range("M8").Value = 20087206
$ws.Range("N8").Value = 23.9
$ws.Range("O8").Value = 6321717

# Row 9: Intermediate school-level
$ws.Range("A9").Value = "Intermediate school-level"
$ws.Range("B9").Value = "All population groups"
$ws.Range("C9").Value = 21127138
$ws.Range("D9").Value = 31
$ws.Range("E9").Value = 6555234
$ws.Range("F9").Value = 0.2
$ws.Range("G9").Value = 33804
$ws.Range("H9").Value = 0.5
$ws.Range("I9").Value = 99615
$ws.Range("J9").Value = 0.9
$ws.Range("K9").Value = 188863
$ws.Range("L9").Value = 67.4
$ws.Range("M9").Value = 14249623
$ws.Range("N9").Value = 32.6
$ws.Range("O9").Value = 6877515

# Row 10: Secondary school
$ws.Range("A10").Value = "Secondary school"
$ws.Range("B10").Value = "All population groups"
$ws.Range("C10").Value = 15845354
$ws.Range("D10").Value = 49.5
$ws.Range("E10").Value = 7851209
$ws.Range("F10").Value = 0.2
$ws.Range("G10").Value = 31752
$ws.Range("H10").Value = 0.5
$ws.Range("I10").Value = 71926
$ws.Range("J10").Value = 0.3
$ws.Range("K10").Value = 44492
$ws.Range("L10").Value = 49.5
$ws.Range("M10").Value = 7845974
$ws.Range("N10").Value = 50.5
$ws.Range("O10").Value = 7999380
